$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I holds student_id / staff_id values for rows 2-7.
# Copy formatting (style) from column B (general text style used for data rows)
# onto the new I cells, then fill in the values.

$ws.Range("B2:B7").Copy()
$ws.Range("I2:I7").PasteSpecial(-4122)

$ws.Range("I2").Value = "18-0246"
$ws.Range("I3").Value = "18-0208"
$ws.Range("I4").Value = "18-0087"
$ws.Range("I5").Value = "18-0054"
$ws.Range("I6").Value = "18-0218 "
$ws.Range("I7").Value = "18-0252"

# Match the selection recorded in the saved workbook after this edit.
$ws.Range("I2:I7").Select()
